$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit reshuffles the per-row data (Fecha=D, Volumen=J, Precio minimo=K,
# Precio maximo=L, Precio promedio ponderado=M, Precio $/Kg=P) across rows 2-9,
# while all other columns stay unchanged.
#
# Row => (D, J, K, L, M, P)
$rows = @{
    2 = @(44200, 1500, 1400, 1500, 1450, 1450)
    3 = @(44210, 1450, 1600, 1700, 1650, 1650)
    4 = @(44175, 1400, 1900, 2000, 1950, 1950)
    5 = @(44895,  200, 1200, 1300, 1255, 1255)
    6 = @(44638,  800, 2500, 2800, 2650, 2650)
    7 = @(44883,  290, 1400, 1500, 1434, 1434)
    8 = @(44893, 3300, 1200, 1300, 1261, 1261)
    9 = @(44537,  800, 1300, 1400, 1350, 1350)
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $ws.Range("D$r").Value = $vals[0]
    $ws.Range("J$r").Value = $vals[1]
    $ws.Range("K$r").Value = $vals[2]
    $ws.Range("L$r").Value = $vals[3]
    $ws.Range("M$r").Value = $vals[4]
    $ws.Range("P$r").Value = $vals[5]
}
